$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "Holcim" vendor above the current row 9 (Petersen Aluminum Corp.),
# pushing everything below it down by one row and keeping the list alphabetical.
$ws.Rows("9").Insert()
$ws.Range("A9").Value = "Holcim"
$ws.Range("G9").Value = "x"

# Update the date/number value in I1
$ws.Range("I1").Value = 33325

# Move the "use this vendor" marks (column B) off the old vendors (rows shifted down by the
# insert: Pro Fastening Systems -> row 11, Stevenson Crane -> row 13, Turek & Sons -> row 16)
$ws.Range("B11").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("B16").ClearContents()

# Turek & Sons no longer needs an envelope
$ws.Range("F16").ClearContents()

# ...and place them on the new vendors instead (ABC Supply, Beacon, Gemco Supply)
$ws.Range("B2").Value = "x"
$ws.Range("B4").Value = "x"
$ws.Range("B8").Value = "x"

# Leave the cursor where the user last clicked
$ws.Range("B8").Select() | Out-Null
